$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = "Normal"
}

Set-TextValue $ws "D2" '60.792.95'
Set-TextValue $ws "E2" '  -0.94%  '
Set-TextValue $ws "D3" '2.905.47'
Set-TextValue $ws "E3" '  -1.81%  '
Set-TextValue $ws "E4" '  +0.00%  '
Set-TextValue $ws "D5" '529.54'
Set-TextValue $ws "E5" '  -1.58%  '
Set-TextValue $ws "D6" '144.83'
Set-TextValue $ws "E6" '  -4.93%  '
Set-TextValue $ws "E7" '  -0.07%  '
Set-TextValue $ws "D8" '0.556'
Set-TextValue $ws "E8" '  -1.15%  '
Set-TextValue $ws "D9" '2.914.03'
Set-TextValue $ws "E9" '  -1.78%  '
Set-TextValue $ws "E10" '  -3.00%  '
Set-TextValue $ws "D11" '6.01'
Set-TextValue $ws "E11" '  -1.72%  '
Set-TextValue $ws "D12" '0.365'
Set-TextValue $ws "E12" '  -0.19%  '
Set-TextValue $ws "D13" '3.413.29'
Set-TextValue $ws "E13" '  -1.98%  '
Set-TextValue $ws "D14" '0.126'
Set-TextValue $ws "E14" '  +1.21%  '
Set-TextValue $ws "D15" '60.776.33'
Set-TextValue $ws "E15" '  -1.02%  '
Set-TextValue $ws "D16" '22.79'
Set-TextValue $ws "E16" '  -3.79%  '
Set-TextValue $ws "D17" '2.910.03'
Set-TextValue $ws "E17" '  -1.73%  '
Set-TextValue $ws "D18" '0.0000142'
Set-TextValue $ws "E18" '  -2.50%  '
Set-TextValue $ws "D19" '5.05'
Set-TextValue $ws "E19" '  -1.29%  '
Set-TextValue $ws "D20" '11.71'
Set-TextValue $ws "E20" '  -1.86%  '
Set-TextValue $ws "D21" '363.33'
Set-TextValue $ws "E21" '  -4.80%  '
Set-TextValue $ws "D22" '6.65'
Set-TextValue $ws "E22" '  +0.12%  '
Set-TextValue $ws "E23" '  -0.06%  '
Set-TextValue $ws "B24" 'LEO'
Set-TextValue $ws "C24" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws "D24" '5.71'
Set-TextValue $ws "E24" '  +0.57%  '
Set-TextValue $ws "B25" 'Litecoin'
Set-TextValue $ws "C25" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws "D25" '64.65'
Set-TextValue $ws "E25" '  -0.40%  '
Set-TextValue $ws "B26" 'Polygon'
Set-TextValue $ws "C26" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws "D26" '0.456'
Set-TextValue $ws "E26" '  -2.59%  '
Set-TextValue $ws "B27" 'Kaspa'
Set-TextValue $ws "C27" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws "D27" '0.183'
Set-TextValue $ws "E27" '  -1.44%  '
Set-TextValue $ws "B28" 'Binance-PegBSC-USD'
Set-TextValue $ws "C28" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws "D28" '1.00'
Set-TextValue $ws "E28" '  +0.00%  '
Set-TextValue $ws "B29" 'InternetComputer(DFINITY)'
Set-TextValue $ws "C29" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws "D29" '7.84'
Set-TextValue $ws "E29" '  -5.31%  '
Set-TextValue $ws "B30" 'PEPE'
Set-TextValue $ws "C30" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws "D30" '0.0₃0867'
Set-TextValue $ws "E30" '  -6.26%  '
Set-TextValue $ws "B31" 'USDe'
Set-TextValue $ws "C31" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws "D31" '0.999'
Set-TextValue $ws "E31" '  +0.00%  '
Set-TextValue $ws "B32" 'PancakeSwap'
Set-TextValue $ws "C32" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws "D32" '1.68'
Set-TextValue $ws "E32" '  -1.41%  '
Set-TextValue $ws "B33" 'EthereumClassic'
Set-TextValue $ws "C33" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws "D33" '19.74'
Set-TextValue $ws "E33" '  -2.61%  '
Set-TextValue $ws "B34" 'Monero'
Set-TextValue $ws "C34" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws "D34" '150.58'
Set-TextValue $ws "E34" '  -3.68%  '
Set-TextValue $ws "B35" 'NEARProtocol'
Set-TextValue $ws "C35" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws "D35" '4.41'
Set-TextValue $ws "E35" '  -3.84%  '
Set-TextValue $ws "B36" 'Aptos'
Set-TextValue $ws "C36" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws "D36" '5.59'
Set-TextValue $ws "E36" '  -6.39%  '
Set-TextValue $ws "B37" 'Fetch.AI'
Set-TextValue $ws "C37" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws "D37" '1.01'
Set-TextValue $ws "E37" '  -4.76%  '
Set-TextValue $ws "B38" 'ImmutableX'
Set-TextValue $ws "C38" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws "D38" '1.21'
Set-TextValue $ws "E38" '  -5.14%  '
Set-TextValue $ws "B39" 'OKB'
Set-TextValue $ws "C39" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws "D39" '37.67'
Set-TextValue $ws "E39" '  +2.48%  '
Set-TextValue $ws "B40" 'Stacks'
Set-TextValue $ws "C40" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws "D40" '1.49'
Set-TextValue $ws "E40" '  -3.05%  '
Set-TextValue $ws "B41" 'Filecoin'
Set-TextValue $ws "C41" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws "D41" '3.74'
Set-TextValue $ws "E41" '  -4.50%  '
Set-TextValue $ws "B42" 'Maker'
Set-TextValue $ws "C42" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws "D42" '2.294.41'
Set-TextValue $ws "E42" '  -5.18%  '
Set-TextValue $ws "B43" 'Mantle'
Set-TextValue $ws "C43" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws "D43" '0.648'
Set-TextValue $ws "E43" '  -1.99%  '
Set-TextValue $ws "B44" 'Hedera'
Set-TextValue $ws "C44" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws "D44" '0.0582'
Set-TextValue $ws "E44" '  -1.77%  '
Set-TextValue $ws "B45" 'EnergySwap'
Set-TextValue $ws "C45" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws "D45" '20.73'
Set-TextValue $ws "E45" '  -6.27%  '
Set-TextValue $ws "B46" 'FirstDigitalUSD'
Set-TextValue $ws "C46" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws "D46" '0.997'
Set-TextValue $ws "E46" '  -0.04%  '
Set-TextValue $ws "B47" 'RenderToken'
Set-TextValue $ws "C47" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws "D47" '5.03'
Set-TextValue $ws "E47" '  +2.75%  '
Set-TextValue $ws "B48" 'VeChain'
Set-TextValue $ws "C48" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws "D48" '0.0238'
Set-TextValue $ws "E48" '  -2.83%  '
Set-TextValue $ws "B49" 'Stellar'
Set-TextValue $ws "C49" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws "D49" '0.0928'
Set-TextValue $ws "E49" '  -2.41%  '
Set-TextValue $ws "B50" 'WhiteBITCoin'
Set-TextValue $ws "C50" 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws "D50" '10.32'
Set-TextValue $ws "E50" '  -1.55%  '
Set-TextValue $ws "B51" 'Bittensor'
Set-TextValue $ws "C51" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws "D51" '252.14'
Set-TextValue $ws "E51" '  -4.11%  '

